$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row on the sheet (data starts at row 2).
$lastRow = 119

# Columns that may hold a HYPERLINK(...) formula whose display-text
# (2nd) argument needs to be added.
$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {

    # 1) "Förändrad" (column C) moves from 45184 to 45186 for every row.
    $ws.Range("C" + $r).Value = 45186

    # 2) Every HYPERLINK formula on the row gets the "Beteckning" (column A)
    #    text added as its second argument, e.g.
    #    HYPERLINK("...", "A 4801-2019")
    $label = $ws.Range("A" + $r).Value()

    foreach ($col in $linkCols) {
        $addr = $col + $r
        $cell = $ws.Range($addr)
        if ($cell.HasFormula()) {
            $f = $cell.Formula()
            if ($f.IndexOf(",") -lt 0) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
